# Updated cryptos list on Sat Jul  1 23:17:43 UTC 2023 with GitHub Actions
#
# Applies the price / percentage-change / coin-identity updates described by
# the diff. Every Price (D) and Volume(1h) (E) cell in this sheet is stored
# as text (inline strings in the original file), so numeric-looking price
# strings are entered with a leading apostrophe to keep Excel from
# re-interpreting them as numbers (which would silently drop things like
# trailing zeros, e.g. "0.4480" -> 0.448).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.592.38"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.922.59"
$ws.Range("E3").Value = "  -0.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'247.31"
$ws.Range("E5").Value = "  +2.76%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4725"
$ws.Range("E7").Value = "  -0.81%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2915"
$ws.Range("E8").Value = "  +1.38%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06844"
$ws.Range("E9").Value = "  +3.04%  "

# Row 10 - Litecoin
$ws.Range("D10").Value = "'105.87"
$ws.Range("E10").Value = "  -2.10%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'18.50"
$ws.Range("E11").Value = "  -3.27%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.929.86"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13 - TRON
$ws.Range("D13").Value = "'0.07734"
$ws.Range("E13").Value = "  +1.62%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.337"
$ws.Range("E14").Value = "  +3.21%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.6728"
$ws.Range("E15").Value = "  +1.60%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "'289.51"
$ws.Range("E16").Value = "  -5.86%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.624.72"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000007651"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19 - Dai
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  -0.07%  "

# Row 20 - was Avalanche, now Uniswap
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'5.574"
$ws.Range("E20").Value = "  +5.25%  "

# Row 21 - was Uniswap, now Avalanche
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'12.97"
$ws.Range("E21").Value = "  -0.41%  "

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "2.180.42"
$ws.Range("E22").Value = "  +0.49%  "

# Row 23 - BinanceUSD
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'6.498"
$ws.Range("E24").Value = "  +3.03%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "'9.544"
$ws.Range("E25").Value = "  +2.28%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'167.05"
$ws.Range("E26").Value = "  -0.54%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'20.82"
$ws.Range("E27").Value = "  +1.74%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "'2.139"
$ws.Range("E28").Value = "  +4.39%  "

# Row 29 - Stellar
$ws.Range("D29").Value = "'0.1071"
$ws.Range("E29").Value = "  -3.06%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "'1.406"
$ws.Range("E30").Value = "  +2.64%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'4.204"
$ws.Range("E31").Value = "  +2.73%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.071"
$ws.Range("E32").Value = "  +3.27%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.05059"
$ws.Range("E33").Value = "  +0.72%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "'0.7354"
$ws.Range("E34").Value = "  -1.01%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.147"
$ws.Range("E35").Value = "  -0.72%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "'0.02061"
$ws.Range("E36").Value = "  +4.96%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "'2.743"
$ws.Range("E37").Value = "  -0.45%  "

# Row 38 - Frax
$ws.Range("D38").Value = "'0.9996"
$ws.Range("E38").Value = "  +0.01%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "'2.686"
$ws.Range("E39").Value = "  -0.23%  "

# Row 40 - Quant
$ws.Range("D40").Value = "'111.81"
$ws.Range("E40").Value = "  +3.71%  "

# Row 41 - RenderToken
$ws.Range("D41").Value = "'2.052"
$ws.Range("E41").Value = "  +0.43%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "'0.4480"
$ws.Range("E42").Value = "  +6.41%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'0.8730"
$ws.Range("E43").Value = "  -0.99%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "'5.912"
$ws.Range("E44").Value = "  +2.09%  "

# Row 45 - PaxDollar
$ws.Range("E45").Value = "  +0.08%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'68.06"
$ws.Range("E46").Value = "  -3.06%  "

# Row 47 - Aptos
$ws.Range("D47").Value = "'7.306"
$ws.Range("E47").Value = "  +0.34%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'9.466"
$ws.Range("E48").Value = "  +2.76%  "

# Row 49 - was Algorand, now BitcoinSV
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'49.28"
$ws.Range("E49").Value = "  +14.97%  "

# Row 50 - was BitcoinSV, now Algorand
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1263"
$ws.Range("E50").Value = "  +3.92%  "

# Row 51 - was Decentraland, now Elrond
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.37"
$ws.Range("E51").Value = "  +1.37%  "
